$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds the "Förändrad" (changed) date serial for every data row (rows 2-533).
# The commit bumps this date by one day (45189 -> 45190) for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 533 }

$ws.Range("C2:C$lastRow").Value = 45190
